$wb = $excel.ActiveWorkbook

# Both the "展览" and "全部类型" sheets contain the same table and need the
# same three cell updates to column F ("想去人数" / number of people interested).
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 493
    $ws.Range("F3").Value = 3351
    $ws.Range("F5").Value = 666
}
